$d = $word.ActiveDocument

# --- 1. Wording tweaks (reviewed with Albert) -------------------------------
# "Register new account with various price ranges" -> "...various account types"
$d.Content.Find.Execute(
    "Register new account with various price ranges", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Register new account with various account types", 2) | Out-Null

# "Price: $10 " -> "Registration price: $10 "
$d.Content.Find.Execute(
    "Price: `$10 ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Registration price: `$10 ", 2) | Out-Null

# "Price: $7" -> "Registration price: $7"
$d.Content.Find.Execute(
    "Price: `$7", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Registration price: `$7", 2) | Out-Null

# "Price: $5" -> "Registration price: $5"
$d.Content.Find.Execute(
    "Price: `$5", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Registration price: `$5", 2) | Out-Null

# --- 2. Swap the "mark books as interesting" bullet for a wish-list bullet --
$d.Content.Find.Execute(
    "Account can mark books as interesting for them, if not available they should be informed when they become available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Account can create wish list by adding books to them", 2) | Out-Null

# --- 3. Move the (hidden) _GoBack bookmark from the end of the "Account is ---
#     able to borrow..." bullet down to the end of the document's last
#     bullet ("Partial searches are allowed...").
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
